# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Swap displayed country names for rows whose totals crossed each other
# - Update the updated statistics (Casos totales, Nuevos casos, Casos activos,
#   Recuperados, Casos criticos, Muertes) for the affected countries

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Septiembre de 2020 a las 22:48"

# --- Country name swaps (rank changed because of updated totals) ---
$ws.Range("A51").Value = "Etiopia"
$ws.Range("A52").Value = "Portugal"

$ws.Range("A56").Value = "Costa Rica"
$ws.Range("A57").Value = "Nigeria"

$ws.Range("A119").Value = "Cabo Verde"
$ws.Range("A120").Value = "Nicaragua"

$ws.Range("A190").Value = "Monaco"
$ws.Range("A191").Value = "Bermudas"

# --- Updated numeric data (row -> B,C,D,E,F,G,H) ---
# Estados Unidos
$ws.Range("B4").Value = 6741429
$ws.Range("C4").Value = 32971
$ws.Range("D4").Value = 4012898
$ws.Range("E4").Value = 2529683
$ws.Range("G4").Value = 328
$ws.Range("H4").Value = 198848

# India
$ws.Range("B5").Value = 4926914
$ws.Range("C5").Value = 81911
$ws.Range("D5").Value = 3856246
$ws.Range("E5").Value = 989860
$ws.Range("G5").Value = 1054
$ws.Range("H5").Value = 80808

# Alemania
$ws.Range("B25").Value = 263210
$ws.Range("C25").Value = 1912
$ws.Range("E25").Value = 18075
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 9435

# Israel
$ws.Range("B27").Value = 160368
$ws.Range("C27").Value = 4764
$ws.Range("D27").Value = 118570
$ws.Range("E27").Value = 40662

# Etiopia (row 51)
$ws.Range("B51").Value = 64786
$ws.Range("C51").Value = 485
$ws.Range("D51").Value = 25333
$ws.Range("E51").Value = 38431
$ws.Range("G51").Value = 9
$ws.Range("H51").Value = 1022

# Portugal (row 52)
$ws.Range("B52").Value = 64596
$ws.Range("C52").Value = 613
$ws.Range("D52").Value = 44185
$ws.Range("E52").Value = 18540
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = 1871

# Costa Rica (row 56)
$ws.Range("B56").Value = 57361
$ws.Range("C56").Value = 937
$ws.Range("D56").Value = 21206
$ws.Range("E56").Value = 35534
$ws.Range("G56").Value = 15
$ws.Range("H56").Value = 621

# Nigeria (row 57)
$ws.Range("B57").Value = 56256
$ws.Range("D57").Value = 44152
$ws.Range("E57").Value = 11022
$ws.Range("H57").Value = 1082

# Costa de Marfil
$ws.Range("B83").Value = 19066
$ws.Range("C83").Value = 53
$ws.Range("D83").Value = 18174
$ws.Range("E83").Value = 772

# Sudan
$ws.Range("B90").Value = 13535
$ws.Range("C90").Value = 19
$ws.Range("D90").Value = 6759
$ws.Range("E90").Value = 5940
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 836

# Guinea
$ws.Range("B96").Value = 10061
$ws.Range("C96").Value = 16
$ws.Range("D96").Value = 9352
$ws.Range("E96").Value = 646

# Mauritania
$ws.Range("B106").Value = 7295
$ws.Range("C106").Value = 19
$ws.Range("D106").Value = 6835
$ws.Range("E106").Value = 299

# Cabo Verde (row 119)
$ws.Range("B119").Value = 4839
$ws.Range("C119").Value = 26
$ws.Range("D119").Value = 4240
$ws.Range("E119").Value = 554
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 45

# Nicaragua (row 120)
$ws.Range("B120").Value = 4818
$ws.Range("D120").Value = 2913
$ws.Range("E120").Value = 1761
$ws.Range("H120").Value = 144

# Ruanda
$ws.Range("B123").Value = 4602
$ws.Range("C123").Value = 11
$ws.Range("D123").Value = 2736
$ws.Range("E123").Value = 1844

# Siria
$ws.Range("B127").Value = 3576
$ws.Range("C127").Value = 36
$ws.Range("D127").Value = 858
$ws.Range("E127").Value = 2561
$ws.Range("G127").Value = 2
$ws.Range("H127").Value = 157

# Aruba
$ws.Range("B139").Value = 3060
$ws.Range("C139").Value = 14
$ws.Range("D139").Value = 1566
$ws.Range("E139").Value = 1474
$ws.Range("G139").Value = 2
$ws.Range("H139").Value = 20

# Togo
$ws.Range("B157").Value = 1578
$ws.Range("C157").Value = 6
$ws.Range("D157").Value = 1204
$ws.Range("E157").Value = 334
$ws.Range("G157").Value = 3
$ws.Range("H157").Value = 40

# Monaco (row 190)
$ws.Range("C190").Value = 6
$ws.Range("D190").Value = 132
$ws.Range("E190").Value = 44
$ws.Range("H190").Value = 1

# Bermudas (row 191)
$ws.Range("B191").Value = 177
$ws.Range("D191").Value = 161
$ws.Range("E191").Value = 7
$ws.Range("H191").Value = 9

# Curazao
$ws.Range("B192").Value = 161
$ws.Range("C192").Value = 4
$ws.Range("D192").Value = 58
$ws.Range("E192").Value = 102
